# Auto-generated Excel COM-interop script to update market-data-driven
# profit calculation columns (H-N) across multiple Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4655.8823
$ws.Range("J62").Value = 4997
$ws.Range("L62").Value = 4997
$ws.Range("N62").Value = -6245
$ws.Range("H64").Value = 5482.75
$ws.Range("I64").Value = 5452
$ws.Range("J64").Value = 5487.143
$ws.Range("K64").Value = 5452
$ws.Range("L64").Value = 5487.143
$ws.Range("M64").Value = -5204
$ws.Range("N64").Value = -5983.143
$ws.Range("H65").Value = 4655.8823
$ws.Range("J65").Value = 4997
$ws.Range("L65").Value = 24985
$ws.Range("N65").Value = -31225
$ws.Range("H67").Value = 5482.75
$ws.Range("I67").Value = 5452
$ws.Range("J67").Value = 5487.143
$ws.Range("K67").Value = 5452
$ws.Range("L67").Value = 5487.143
$ws.Range("M67").Value = -4594
$ws.Range("N67").Value = -7203.143
$ws.Range("H81").Value = 195000
$ws.Range("J81").Value = 195000
$ws.Range("L81").Value = 195000
$ws.Range("N81").Value = -196996
$ws.Range("H84").Value = 195000
$ws.Range("J84").Value = 195000
$ws.Range("L84").Value = 585000
$ws.Range("N84").Value = -594984
$ws.Range("H116").Value = 3596.125
$ws.Range("J116").Value = 3890.5
$ws.Range("L116").Value = 3890.5
$ws.Range("N116").Value = -10774.5
$ws.Range("H125").Value = 7100
$ws.Range("J125").Value = 7625
$ws.Range("L125").Value = 68625
$ws.Range("N125").Value = -73545
$ws.Range("H132").Value = 4973.3794
$ws.Range("I132").Value = 4474.6
$ws.Range("J132").Value = 8090.75
$ws.Range("K132").Value = 13423.8
$ws.Range("L132").Value = 24272.25
$ws.Range("M132").Value = -10893.8
$ws.Range("N132").Value = -29332.25
$ws.Range("H135").Value = 438.69232
$ws.Range("I135").Value = 326.25
$ws.Range("J135").Value = 618.6
$ws.Range("K135").Value = 2936.25
$ws.Range("L135").Value = 5567.400000000001
$ws.Range("M135").Value = -401.25
$ws.Range("N135").Value = -10637.4
$ws.Range("H137").Value = 1702.7916
$ws.Range("I137").Value = 1650.7826
$ws.Range("J137").Value = 2899
$ws.Range("K137").Value = 4952.3478
$ws.Range("L137").Value = 8697
$ws.Range("M137").Value = -2402.3478
$ws.Range("N137").Value = -13797
$ws.Range("H138").Value = 3900.3157
$ws.Range("I138").Value = 3686.375
$ws.Range("J138").Value = 3957.3667
$ws.Range("K138").Value = 11059.125
$ws.Range("L138").Value = 11872.1001
$ws.Range("M138").Value = -5919.125
$ws.Range("N138").Value = -22152.1001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17761.422
$ws.Range("I32").Value = 21151.326
$ws.Range("J32").Value = 6687.7334
$ws.Range("K32").Value = 21151.326
$ws.Range("L32").Value = 6687.7334
$ws.Range("M32").Value = -20864.326
$ws.Range("N32").Value = -7261.7334
$ws.Range("H45").Value = 3809.75
$ws.Range("I45").Value = 2348.9285
$ws.Range("K45").Value = 2348.9285
$ws.Range("M45").Value = -1971.9285
$ws.Range("H97").Value = 1266.0769
$ws.Range("I97").Value = 922.375
$ws.Range("J97").Value = 1816
$ws.Range("K97").Value = 922.375
$ws.Range("L97").Value = 1816
$ws.Range("M97").Value = -426.375
$ws.Range("N97").Value = -2808
$ws.Range("H132").Value = 58728.527
$ws.Range("I132").Value = 205599.8
$ws.Range("J132").Value = 6274.5
$ws.Range("K132").Value = 616799.3999999999
$ws.Range("L132").Value = 18823.5
$ws.Range("M132").Value = -614269.3999999999
$ws.Range("N132").Value = -23883.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 39998.25
$ws.Range("J6").Value = 39998.25
$ws.Range("L6").Value = 39998.25
$ws.Range("N6").Value = -40224.25
$ws.Range("H20").Value = 2613.6
$ws.Range("I20").Value = 1071.5
$ws.Range("J20").Value = 5697.8
$ws.Range("K20").Value = 1071.5
$ws.Range("L20").Value = 5697.8
$ws.Range("M20").Value = -824.5
$ws.Range("N20").Value = -6191.8
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H134").Value = 3358
$ws.Range("I134").Value = 3358
$ws.Range("K134").Value = 10074
$ws.Range("M134").Value = -7539

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1771.7587
$ws.Range("I31").Value = 1813.25
$ws.Range("J31").Value = 1572.6
$ws.Range("K31").Value = 1813.25
$ws.Range("L31").Value = 1572.6
$ws.Range("M31").Value = -1518.25
$ws.Range("N31").Value = -2162.6
$ws.Range("H34").Value = 1771.7587
$ws.Range("I34").Value = 1813.25
$ws.Range("J34").Value = 1572.6
$ws.Range("K34").Value = 1813.25
$ws.Range("L34").Value = 1572.6
$ws.Range("M34").Value = -1611.25
$ws.Range("N34").Value = -1976.6
$ws.Range("H134").Value = 29403.361
$ws.Range("I134").Value = 33694.547
$ws.Range("K134").Value = 101083.641
$ws.Range("M134").Value = -98548.641

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 857.75
$ws.Range("I5").Value = 439.4
$ws.Range("J5").Value = 1047.909
$ws.Range("K5").Value = 1318.2
$ws.Range("L5").Value = 3143.727
$ws.Range("M5").Value = -1206.2
$ws.Range("N5").Value = -3367.727
$ws.Range("H113").Value = 2107.7058
$ws.Range("I113").Value = 824.73914
$ws.Range("J113").Value = 4790.273
$ws.Range("K113").Value = 2474.21742
$ws.Range("L113").Value = 14370.819
$ws.Range("M113").Value = -304.2174199999999
$ws.Range("N113").Value = -18710.819
$ws.Range("H134").Value = 1833.5333
$ws.Range("I134").Value = 1250.3
$ws.Range("K134").Value = 3750.9
$ws.Range("M134").Value = 1319.1
$ws.Range("H135").Value = 857.75
$ws.Range("I135").Value = 439.4
$ws.Range("J135").Value = 1047.909
$ws.Range("K135").Value = 3954.6
$ws.Range("L135").Value = 9431.181
$ws.Range("M135").Value = -1419.6
$ws.Range("N135").Value = -14501.181
$ws.Range("H140").Value = 4523.3335
$ws.Range("I140").Value = 4281.7
$ws.Range("K140").Value = 12845.1
$ws.Range("M140").Value = -7665.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 35998
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H97").Value = 636
$ws.Range("J97").Value = 935.5
$ws.Range("L97").Value = 935.5
$ws.Range("N97").Value = -1927.5
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H126").Value = 5141.125
$ws.Range("I126").Value = 4995.8
$ws.Range("K126").Value = 14987.4
$ws.Range("M126").Value = -12517.4
$ws.Range("H133").Value = 130780
$ws.Range("J133").Value = 130780
$ws.Range("L133").Value = 130780
$ws.Range("N133").Value = -140900
$ws.Range("H140").Value = 89998.664
$ws.Range("J140").Value = 89998.664
$ws.Range("L140").Value = 89998.664
$ws.Range("N140").Value = -100358.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5574.5
$ws.Range("H48").Value = 27805.6
$ws.Range("J48").Value = 36666
$ws.Range("L48").Value = 36666
$ws.Range("N48").Value = -37988
$ws.Range("H55").Value = 533.8182
$ws.Range("I55").Value = 359.6
$ws.Range("K55").Value = 359.6
$ws.Range("M55").Value = -186.6
$ws.Range("H106").Value = 19997
$ws.Range("J106").Value = 19997
$ws.Range("L106").Value = 19997
$ws.Range("N106").Value = -22521
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H126").Value = 5574.5
$ws.Range("H132").Value = 54260.652
$ws.Range("I132").Value = 56479.59
$ws.Range("K132").Value = 169438.77
$ws.Range("M132").Value = -166908.77
$ws.Range("H136").Value = 5263.7144
$ws.Range("I136").Value = 4391
$ws.Range("J136").Value = 10500
$ws.Range("K136").Value = 13173
$ws.Range("L136").Value = 31500
$ws.Range("M136").Value = -10623
$ws.Range("N136").Value = -36600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 31123.777
$ws.Range("J111").Value = 31123.777
$ws.Range("L111").Value = 31123.777
$ws.Range("N111").Value = -39303.777
$ws.Range("H136").Value = 2327.8076
$ws.Range("I136").Value = 1549.3043
$ws.Range("J136").Value = 8296.333000000001
$ws.Range("K136").Value = 4647.9129
$ws.Range("L136").Value = 24888.999
$ws.Range("M136").Value = -2097.9129
$ws.Range("N136").Value = -29988.999

